$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in columns D and E hold numeric-looking text (prices / percentages) that
# must stay plain text (matching the source inline-string cells), so force the
# number format to Text before writing the new literal value.
$textCells = @{
    "D2" = "310.21"
    "E2" = "1.34%"
    "D3" = "35.54"
    "E3" = "-1.85%"
    "D4" = "5.122"
    "E4" = "1.41%"
    "D5" = "0.08196"
    "E5" = "3.64%"
    "D6" = "2.045"
    "E6" = "-10.78%"
    "D7" = "7.965"
    "E7" = "-0.36%"
    "D8" = "4.131"
    "E8" = "-0.47%"
    "D9" = "2.912"
    "E9" = "10.66%"
    "D10" = "0.9277"
    "E10" = "-0.09%"
    "D11" = "0.1083"
    "E11" = "9.89%"
    "D12" = "0.1922"
    "E12" = "3.09%"
    "D13" = "0.09439"
    "E13" = "4.93%"
    "D14" = "0.03607"
    "E14" = "-3.97%"
    "D15" = "0.09885"
    "E15" = "-0.39%"
    "D16" = "0.001431"
    "E16" = "-0.01%"
    "D17" = "0.005719"
    "E17" = "0.43%"
    "D18" = "3.468"
    "E18" = "1.00%"
    "D19" = "0.3417"
    "E19" = "1.46%"
    "D20" = "0.1313"
    "E20" = "-0.52%"
    "D21" = "5.102"
    "E21" = "0.64%"
    "D22" = "0.2192"
    "E22" = "-2.62%"
    "D23" = "0.04551"
    "E23" = "-0.57%"
    "D24" = "0.001225"
    "E24" = "-0.78%"
    "D25" = "0.004784"
    "E25" = "0.21%"
    "D26" = "0.0001251"
    "E26" = "-3.76%"
    "D27" = "0.0004452"
    "E27" = "-6.05%"
    "D39" = "0.01972"
    "E39" = "2.92%"
    "D40" = "0.04897"
    "E40" = "-0.39%"
    "D41" = "0.007821"
    "E41" = "-0.62%"
    "D42" = "0.009811"
    "E42" = "25.58%"
    "E43" = "-0.85%"
    "D44" = "0.002117"
    "E44" = "-0.09%"
    "E45" = "1.14%"
    "D46" = "0.00006504"
    "E46" = "6.01%"
    "D47" = "0.00000000751"
    "E47" = "-0.04%"
    "D48" = "64.60"
    "E48" = "24.80%"
    "E49" = "-16.83%"
    "D50" = "0.00002102"
    "E50" = "-0.04%"
    "D51" = "0.0002002"
    "E51" = "-0.04%"
}

foreach ($addr in $textCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $textCells[$addr]
}

# Columns B and C are plain (non-numeric-looking) text, Excel keeps these as
# strings natively so no special handling is required.
$textValues = @{
    "B8" = "GateToken"
    "C8" = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
    "B9" = "BTSEToken"
    "C9" = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
    "B10" = "MXToken"
    "C10" = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
    "B11" = "LiechtensteinCryptoassetsExchange"
    "C11" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "B12" = "WazirX"
    "C12" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "B13" = "MandalaExchangeToken"
    "C13" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "B14" = "BitrueCoin"
    "C14" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "B15" = "BitMartToken"
    "C15" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "B16" = "BitForexToken"
    "C16" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "B17" = "TigerCash"
    "C17" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "B18" = "LEO"
    "C18" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
}

foreach ($addr in $textValues.Keys) {
    $ws.Range($addr).Value = $textValues[$addr]
}
